# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计",
#    populated with the fund-holding detail for the new quarter.
# 2. Insert a new summary row at the top of "总计" for "2022-Q1" and
#    renumber the existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" worksheet right after "2021-Q4"
# ---------------------------------------------------------------------
$prev = $wb.Worksheets.Item("2021-Q4")
$qtr  = $wb.Worksheets.Add($null, $prev)
$qtr.Name = "2022-Q1"

# Borrow the exact header / index-column formatting used by the other
# quarterly sheets (bold, centered, thin-bordered) so no new cell
# styles get introduced.
$prev.Range("B1:H1").Copy()
$qtr.Range("B1:H1").PasteSpecial(-4122)
$prev.Range("A2").Copy()
$qtr.Range("A2:A9").PasteSpecial(-4122)

$qtr.Range("B1").Value = "基金代码"
$qtr.Range("C1").Value = "基金名称"
$qtr.Range("D1").Value = "基金规模"
$qtr.Range("E1").Value = "股票总仓位"
$qtr.Range("F1").Value = "仓位占比"
$qtr.Range("G1").Value = "持有市值(亿元)"
$qtr.Range("H1").Value = "仓位排名"

# index, fund code, fund name, fund scale, total position, position pct,
# held value (亿元), position rank
$rows = @(
    @(0, "001838", "国投瑞银国家安全灵活配置混合", "32.13", "94.68", "7.49", "2.4065", 6),
    @(1, "006440", "中信建投中证500指数增强A",     "5.78",  "94.71", "0.97", "0.0561", 8),
    @(2, "690001", "民生加银品牌蓝筹混合",           "1.21",  "93.01", "3.78", "0.0457", 9),
    @(3, "000714", "诺安稳健回报灵活配置混合A",     "1.96",  "64.55", "2.30", "0.0451", 9),
    @(4, "002052", "诺安稳健回报灵活配置混合C",     "1.63",  "64.55", "2.30", "0.0375", 9),
    @(5, "011685", "创金合信先进装备股票A",         "0.73",  "92.01", "4.39", "0.0320", 9),
    @(6, "006441", "中信建投中证500指数增强C",     "3.11",  "94.71", "0.97", "0.0302", 8),
    @(7, "011686", "创金合信先进装备股票C",         "0.17",  "92.01", "4.39", "0.0075", 9)
)

$r = 2
foreach ($item in $rows) {
    $qtr.Range("A$r").Value = $item[0]

    $qtr.Range("B$r").Value = "'" + $item[1]
    $qtr.Range("B$r").ClearFormats()

    $qtr.Range("C$r").Value = $item[2]

    $qtr.Range("D$r").Value = "'" + $item[3]
    $qtr.Range("D$r").ClearFormats()

    $qtr.Range("E$r").Value = "'" + $item[4]
    $qtr.Range("E$r").ClearFormats()

    $qtr.Range("F$r").Value = "'" + $item[5]
    $qtr.Range("F$r").ClearFormats()

    $qtr.Range("G$r").Value = "'" + $item[6]
    $qtr.Range("G$r").ClearFormats()

    $qtr.Range("H$r").Value = $item[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: prepend a "2022-Q1" row to the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Re-apply the index-column style (bold/border) that the original A2
# cell (now shifted to A3) carries, so the new A2 matches its siblings.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 2.66

# Renumber the index column for the rows that shifted down one place.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
